$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Hunk 1: "...CaseStudy2-data.csv on AWS S3) to..." ->
#         "...CaseStudy2-data.csv on AWS S3 in the smuddsproject2 bucket) to..."
# The appended text carries the same Bold run formatting as " on AWS S3".
# ---------------------------------------------------------------------------
$r1 = $d.Content
$r1.Find.Execute("-data.csv on AWS S3", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r1.Collapse(0)
$r1.InsertAfter(" in the smuddsproject2 bucket")
$r1.Font.Bold = 1

# ---------------------------------------------------------------------------
# Hunk 2: first "An example submission file can be found on GitHub: "
# (precedes "Case2PredictionsClassifyEXAMPLE.csv") ->
# "An example submission file can be found on AWS S3 in the smuddsproject2 bucket: "
# ---------------------------------------------------------------------------
$r2 = $d.Content
$r2.Find.Execute("An example submission file can be found on GitHub: ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r2.Text = "An example submission file can be found "
$r2.Collapse(0)
$r2.InsertAfter("on ")
$r2.Collapse(0)
$r2.InsertAfter("AWS S3 in the smuddsproject2 bucket")
$r2.Collapse(0)
$r2.InsertAfter(": ")
$r2.Collapse(0)

# ---------------------------------------------------------------------------
# Hunk 3: second "An example submission file can be found on GitHub: "
# (precedes "Case2PredictionsRegressEXAMPLE.csv") ->
# "An example submission file can be found on AWS S3 in the smuddsproject 2 bucket: "
# ---------------------------------------------------------------------------
$r3 = $d.Range($r2.End, $d.Content.End)
$r3.Find.Execute("An example submission file can be found on GitHub: ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r3.Text = "An example submission file can be found on "
$r3.Collapse(0)
$r3.InsertAfter("AWS S3 in the ")
$r3.Collapse(0)
$r3.InsertAfter("smuddsproject")
$r3.Collapse(0)
$r3.InsertAfter(" 2 bucket")
$r3.Collapse(0)
$r3.InsertAfter(": ")

# ---------------------------------------------------------------------------
# Hunk 4: "with an RMarkdown" -> "with a RMarkdown"
# ---------------------------------------------------------------------------
$r4 = $d.Content
$r4.Find.Execute("with an ", $true, $false, $false, $false, $false, $true, 1, $false, "with a ", 2) | Out-Null

# ---------------------------------------------------------------------------
# Hunk 5: a <w:lastRenderedPageBreak/> marker is added right before the
# "Due Dates: " run. This is a purely-cosmetic, Word-internal layout artifact
# (the recorded location of the last page break during rendering) that Word
# itself recomputes on every repagination; it has no insertion API in the
# Word object model (COM/VBA), so it cannot be produced from automation code
# here either - it is intentionally left alone.
# ---------------------------------------------------------------------------
